# fix scaling of 2024 logs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1  = 19
    2  = 35
    3  = 44
    4  = 58
    5  = 68
    6  = 72
    7  = 87
    8  = 106
    9  = 112
    10 = 132
    11 = 164
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
